$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "96.036.45"
Set-TextValue "E2" "  -1.80%  "
Set-TextValue "D3" "3.304.82"
Set-TextValue "E3" "  -3.29%  "
Set-TextValue "E4" "  +0.26%  "
Set-TextValue "D5" "245.89"
Set-TextValue "E5" "  -4.08%  "
Set-TextValue "D6" "645.50"
Set-TextValue "E6" "  -1.47%  "
Set-TextValue "D7" "1.34"
Set-TextValue "E7" "  -8.82%  "
Set-TextValue "D8" "0.409"
Set-TextValue "E8" "  -5.02%  "
Set-TextValue "E9" "  +0.20%  "
Set-TextValue "D10" "0.963"
Set-TextValue "E10" "  -8.07%  "
Set-TextValue "D11" "3.303.20"
Set-TextValue "E11" "  -3.32%  "
Set-TextValue "D12" "0.204"
Set-TextValue "E12" "  -4.50%  "
Set-TextValue "D13" "39.07"
Set-TextValue "E13" "  -7.85%  "
Set-TextValue "B14" "WrappedBTC"
Set-TextValue "C14" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D14" "95.823.62"
Set-TextValue "E14" "  -1.87%  "
Set-TextValue "B15" "Toncoin"
Set-TextValue "C15" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D15" "6.05"
Set-TextValue "E15" "  -7.46%  "
Set-TextValue "D16" "0.0000246"
Set-TextValue "E16" "  -4.71%  "
Set-TextValue "D17" "3.915.12"
Set-TextValue "E17" "  -3.04%  "
Set-TextValue "D18" "8.34"
Set-TextValue "E18" "  -2.58%  "
Set-TextValue "D19" "3.296.84"
Set-TextValue "E19" "  -3.20%  "
Set-TextValue "D20" "16.51"
Set-TextValue "E20" "  -6.00%  "
Set-TextValue "D21" "499.79"
Set-TextValue "E21" "  -2.07%  "
Set-TextValue "D22" "0.469"
Set-TextValue "E22" "  -8.65%  "
Set-TextValue "D23" "10.28"
Set-TextValue "E23" "  -5.52%  "
Set-TextValue "D24" "3.26"
Set-TextValue "E24" "  -5.70%  "
Set-TextValue "D25" "0.0000193"
Set-TextValue "E25" "  -6.27%  "
Set-TextValue "D26" "6.31"
Set-TextValue "E26" "  +3.75%  "
Set-TextValue "D27" "93.60"
Set-TextValue "E27" "  -2.82%  "
Set-TextValue "D28" "11.72"
Set-TextValue "E28" "  -7.54%  "
Set-TextValue "D29" "3.490.90"
Set-TextValue "E29" "  -2.05%  "
Set-TextValue "E30" "  +0.30%  "
Set-TextValue "D31" "0.139"
Set-TextValue "E31" "  -8.10%  "
Set-TextValue "D32" "10.57"
Set-TextValue "E32" "  -7.42%  "
Set-TextValue "D33" "0.182"
Set-TextValue "E33" "  -7.08%  "
Set-TextValue "D34" "2.43"
Set-TextValue "E34" "  +10.91%  "
Set-TextValue "D35" "0.998"
Set-TextValue "E35" "  -0.06%  "
Set-TextValue "D36" "0.531"
Set-TextValue "E36" "  -7.05%  "
Set-TextValue "D37" "27.52"
Set-TextValue "E37" "  -7.45%  "
Set-TextValue "D38" "1.47"
Set-TextValue "E38" "  +5.38%  "
Set-TextValue "B39" "USDe"
Set-TextValue "C39" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D39" "1.00"
Set-TextValue "E39" "  +0.00%  "
Set-TextValue "B40" "RenderToken"
Set-TextValue "C40" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D40" "7.37"
Set-TextValue "E40" "  -5.77%  "
Set-TextValue "D41" "0.148"
Set-TextValue "E41" "  -5.34%  "
Set-TextValue "D42" "495.43"
Set-TextValue "E42" "  -3.56%  "
Set-TextValue "D43" "24.43"
Set-TextValue "E43" "  -1.13%  "
Set-TextValue "D44" "3.63"
Set-TextValue "E44" "  -0.78%  "
Set-TextValue "D45" "0.809"
Set-TextValue "E45" "  -4.91%  "
Set-TextValue "D46" "0.0399"
Set-TextValue "E46" "  -5.69%  "
Set-TextValue "D47" "5.31"
Set-TextValue "E47" "  -1.73%  "
Set-TextValue "D48" "8.13"
Set-TextValue "E48" "  -0.61%  "
Set-TextValue "B49" "OKB"
Set-TextValue "C49" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D49" "52.85"
Set-TextValue "E49" "  +4.40%  "
Set-TextValue "B50" "ImmutableX"
Set-TextValue "C50" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D50" "1.58"
Set-TextValue "E50" "  -0.35%  "
Set-TextValue "D51" "3.06"
Set-TextValue "E51" "  -7.55%  "
